$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.872.10"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.375.67"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.62"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.53"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -3.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.44"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -3.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.797.64"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.833.21"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.377.00"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.18"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.14"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.45"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +1.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0760"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.97"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +11.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.401"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.12"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -2.57%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.40"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.59"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "144.65"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +3.32%  "

$ws.Range("E43").Value = "  -3.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0967"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.65"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0511"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("E51").Value = "  -1.42%  "
